$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the two name cells (these become shared strings: "Eihab" -> A1, "Ahmed" -> C4)
$ws.Range("A1").Value = "Eihab"
$ws.Range("C4").Value = "Ahmed"

# Match the saved selection/active cell from the authored workbook (P15)
$ws.Range("P15").Select()
